$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of row -> new dSF (column F) value, per the repulled/recalculated data
$changes = @{
    6  = -4
    7  = -4
    18 = -7
    22 = -10
    24 = -1
    27 = 0
    28 = 6
    30 = 1
    31 = 2
    34 = 4
    39 = -2
    45 = 0
    49 = -2
    53 = 1
    58 = -3
    68 = -1
    69 = 0
    74 = -1
}

foreach ($row in $changes.Keys) {
    $ws.Range("F$row").Value = $changes[$row]
}
